$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 0
$ws.Range("F6").Value = -6
$ws.Range("F11").Value = -3
$ws.Range("F17").Value = 0
$ws.Range("F20").Value = -4
$ws.Range("F21").Value = -4
$ws.Range("F22").Value = -5
$ws.Range("F30").Value = 1
$ws.Range("F33").Value = -1
$ws.Range("F37").Value = -3
$ws.Range("F41").Value = 2
$ws.Range("F42").Value = 2
$ws.Range("F43").Value = 2
$ws.Range("F44").Value = -1
$ws.Range("F45").Value = -1
$ws.Range("F56").Value = -1
$ws.Range("F57").Value = 0
$ws.Range("F63").Value = -4
